# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@80fa500adfae01c9a5dd7ef65e90accc96781b5c 🚀
#
# The StructureDefinition IG spreadsheet was regenerated under the new
# LinuxForHealth org/branding: URL, Version, Date and Publisher metadata
# are bumped, and a stray duplicated FHIR invariant ("ele-1/ext-1") that
# had been erroneously copied onto the root "Extension" element row is
# cleared (it correctly stays on the "Extension.extension" row).

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------
# URL
$wsMetadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/process-record-id"
# Version
$wsMetadata.Range("B3").Value = "8.0.0"
# Date
$wsMetadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
# Publisher
$wsMetadata.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
# The "Fixed Value" for Extension.url mirrors the StructureDefinition URL
# (it shared the same underlying string as Metadata!B2 in the source
# template), so it needs to be updated to match too.
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/process-record-id"

# The "Constraint(s)" cell on the root Extension row (AI2) incorrectly
# duplicated the ele-1/ext-1 invariant text that belongs to the
# Extension.extension row (AI4); clear it.
$wsElements.Range("AI2").Value = ""
